# Add 2022-Q3 data
# 1) Insert a new summary row in the "总计" (Total) sheet for 2022-Q3.
# 2) Insert a new worksheet "2022-Q3" (positioned right after "总计", before
#    "2022-Q2") with the per-fund holding detail for that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" sheet: insert row 2 with the 2022-Q3 totals, shifting the
#    existing quarters down by one row.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# Clear the quote-prefix text formatting Insert() leaves behind on the
# pushed-down blank row, and re-apply the same style as column A uses
# elsewhere in the table (bold/border, from the row that used to be row 2).
$total.Range("B2:D2").ClearFormats()
$total.Cells.Item(3, 1).Copy($total.Cells.Item(2, 1))

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 9
$total.Cells.Item(2, 4).Value = 1.62

# ---------------------------------------------------------------------
# 2) New "2022-Q3" worksheet with the fund-level detail table. Copy the
#    structurally-identical "2022-Q2" sheet's range first so every cell
#    picks up matching styles (header row, column A, numeric rank column),
#    then overwrite with the 2022-Q3 figures.
# ---------------------------------------------------------------------
$existingQ2 = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($existingQ2)
$newSheet.Name = "2022-Q3"

$existingQ2.Range("B1:H10").Copy($newSheet.Range("B1:H10"))
$existingQ2.Range("A2:A10").Copy($newSheet.Range("A2:A10"))

$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1

$rows = @(
    @("001908", "国投瑞银境煊灵活配置混合C", "4.56", "92.52", "8.00", "0.3648", 6),
    @("010425", "国投瑞银开放视角精选混合A", "4.27", "91.65", "7.77", "0.3318", 8),
    @("010673", "兴全中证800六个月持有期指数增强A", "12.02", "96.73", "2.59", "0.3113", 8),
    @("001907", "国投瑞银境煊灵活配置混合A", "3.76", "92.52", "8.00", "0.3008", 6),
    @("015309", "国投瑞银境煊灵活配置混合E", "2.68", "92.52", "8.00", "0.2144", 6),
    @("010426", "国投瑞银开放视角精选混合C", "0.59", "91.65", "7.77", "0.0458", 8),
    @("010674", "兴全中证800六个月持有期指数增强C", "1.29", "96.73", "2.59", "0.0334", 8),
    @("015148", "华安中证1000指数增强A", "1.50", "91.03", "0.67", "0.0100", 6),
    @("015149", "华安中证1000指数增强C", "0.77", "91.03", "0.67", "0.0052", 6)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $newSheet.Cells.Item($r, 2).Value = "'" + $data[0]
    $newSheet.Cells.Item($r, 2).Style = "Normal"

    $newSheet.Cells.Item($r, 3).Value = $data[1]

    $newSheet.Cells.Item($r, 4).Value = "'" + $data[2]
    $newSheet.Cells.Item($r, 4).Style = "Normal"

    $newSheet.Cells.Item($r, 5).Value = "'" + $data[3]
    $newSheet.Cells.Item($r, 5).Style = "Normal"

    $newSheet.Cells.Item($r, 6).Value = "'" + $data[4]
    $newSheet.Cells.Item($r, 6).Style = "Normal"

    $newSheet.Cells.Item($r, 7).Value = "'" + $data[5]
    $newSheet.Cells.Item($r, 7).Style = "Normal"

    $newSheet.Cells.Item($r, 8).Value = $data[6]
}

# Restore the originally active sheet/selection.
$total.Activate()
$total.Range("A1").Select()
